$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add uncertainty formulas in column I for rows 2-6
$ws.Range("I2").Formula = "=SQRT((C2/B2)^2+(E2/D2)^2+(G2/F2)^2)*H2"
$ws.Range("I3:I6").Formula = "=SQRT((C3/B3)^2+(E3/D3)^2+(G3/F3)^2)*H3"

# Update the selected cell/range on the sheet view
$ws.Range("J11").Select()
